$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell ref -> new value, taken verbatim from the target diff.
$updates = @(
    @{ Ref = 'D2'; Value = '26.608.32' },
    @{ Ref = 'E2'; Value = '  -2.17%  ' },
    @{ Ref = 'D3'; Value = '1.585.08' },
    @{ Ref = 'E3'; Value = '  -2.78%  ' },
    @{ Ref = 'E4'; Value = '  -0.04%  ' },
    @{ Ref = 'D5'; Value = '210.96' },
    @{ Ref = 'E5'; Value = '  -2.29%  ' },
    @{ Ref = 'E6'; Value = '  -2.15%  ' },
    @{ Ref = 'E7'; Value = '  -0.03%  ' },
    @{ Ref = 'E8'; Value = '  -2.48%  ' },
    @{ Ref = 'E9'; Value = '  -0.97%  ' },
    @{ Ref = 'D10'; Value = '19.56' },
    @{ Ref = 'E10'; Value = '  -3.30%  ' },
    @{ Ref = 'D11'; Value = '0.0834' },
    @{ Ref = 'E11'; Value = '  -1.66%  ' },
    @{ Ref = 'D12'; Value = '1.805.20' },
    @{ Ref = 'E12'; Value = '  -2.85%  ' },
    @{ Ref = 'D13'; Value = '1.588.19' },
    @{ Ref = 'E13'; Value = '  -2.54%  ' },
    @{ Ref = 'E14'; Value = '  -1.73%  ' },
    @{ Ref = 'E15'; Value = '  -2.87%  ' },
    @{ Ref = 'D16'; Value = '64.47' },
    @{ Ref = 'E16'; Value = '  -0.35%  ' },
    @{ Ref = 'D17'; Value = '26.622.48' },
    @{ Ref = 'E17'; Value = '  -2.02%  ' },
    @{ Ref = 'E18'; Value = '  -0.48%  ' },
    @{ Ref = 'E19'; Value = '  +0.00%  ' },
    @{ Ref = 'D20'; Value = '207.67' },
    @{ Ref = 'E20'; Value = '  -3.89%  ' },
    @{ Ref = 'D21'; Value = '6.77' },
    @{ Ref = 'E21'; Value = '  -2.27%  ' },
    @{ Ref = 'E22'; Value = '  -3.17%  ' },
    @{ Ref = 'D23'; Value = '2.37' },
    @{ Ref = 'E23'; Value = '  -4.61%  ' },
    @{ Ref = 'E24'; Value = '  -2.27%  ' },
    @{ Ref = 'D25'; Value = '146.67' },
    @{ Ref = 'E25'; Value = '  -0.89%  ' },
    @{ Ref = 'E26'; Value = '  -0.04%  ' },
    @{ Ref = 'E27'; Value = '  +1.71%  ' },
    @{ Ref = 'E28'; Value = '  -4.13%  ' },
    @{ Ref = 'D29'; Value = '15.29' },
    @{ Ref = 'E29'; Value = '  -1.92%  ' },
    @{ Ref = 'E30'; Value = '  -0.43%  ' },
    @{ Ref = 'E31'; Value = '  -1.92%  ' },
    @{ Ref = 'E32'; Value = '  -4.11%  ' },
    @{ Ref = 'D33'; Value = '0.680' },
    @{ Ref = 'E33'; Value = '  +25.32%  ' },
    @{ Ref = 'D34'; Value = '2.93' },
    @{ Ref = 'E34'; Value = '  -2.77%  ' },
    @{ Ref = 'D35'; Value = '1.324.04' },
    @{ Ref = 'E35'; Value = '  +0.74%  ' },
    @{ Ref = 'B36'; Value = 'LidoDAOToken' },
    @{ Ref = 'C36'; Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo' },
    @{ Ref = 'D36'; Value = '1.50' },
    @{ Ref = 'E36'; Value = '  -3.56%  ' },
    @{ Ref = 'B37'; Value = 'HuobiToken' },
    @{ Ref = 'C37'; Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht' },
    @{ Ref = 'D37'; Value = '2.42' },
    @{ Ref = 'E37'; Value = '  -1.08%  ' },
    @{ Ref = 'E38'; Value = '  -1.31%  ' },
    @{ Ref = 'D39'; Value = '0.824' },
    @{ Ref = 'E39'; Value = '  -3.08%  ' },
    @{ Ref = 'E40'; Value = '  +0.01%  ' },
    @{ Ref = 'D41'; Value = '5.35' },
    @{ Ref = 'E41'; Value = '  +3.07%  ' },
    @{ Ref = 'E42'; Value = '  -2.24%  ' },
    @{ Ref = 'E43'; Value = '  -3.65%  ' },
    @{ Ref = 'D44'; Value = '63.49' },
    @{ Ref = 'E44'; Value = '  -0.14%  ' },
    @{ Ref = 'D45'; Value = '1.718.74' },
    @{ Ref = 'E45'; Value = '  -2.69%  ' },
    @{ Ref = 'D46'; Value = '89.62' },
    @{ Ref = 'E46'; Value = '  -1.16%  ' },
    @{ Ref = 'E47'; Value = '  +1.08%  ' },
    @{ Ref = 'E48'; Value = '  +3.20%  ' },
    @{ Ref = 'D49'; Value = '0.0988' },
    @{ Ref = 'E49'; Value = '  +3.13%  ' },
    @{ Ref = 'D50'; Value = '0.0506' },
    @{ Ref = 'E50'; Value = '  -1.93%  ' },
    @{ Ref = 'D51'; Value = '7.48' },
    @{ Ref = 'E51'; Value = '  -0.67%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Ref)
    $newValue = $u.Value

    # All touched cells hold plain text in the source workbook (inline/shared
    # strings), even ones that look like plain numbers (e.g. "0.680", "1.50").
    # Assigning a numeric-looking string straight to .Value lets Excel coerce it
    # to a real number (dropping e.g. a trailing zero), so force Text format
    # first for those, then drop back to the default style so no stray
    # per-cell formatting is introduced.
    if ($newValue -match '^[+-]?\d+(\.\d+)?$') {
        $cell.NumberFormat = "@"
        $cell.Value = $newValue
        $cell.Style = "Normal"
    } else {
        $cell.Value = $newValue
    }
}
